# 화면 수정 2025-09-22 2시 20분
# Remove the "닫기"(close) control group from slide 1:
#   - 직사각형 26 (id=27)
#   - 직사각형 30 (id=31)
#   - 그래픽 31 "닫기" picture (id=32)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$idsToDelete = @(27, 31, 32)

foreach ($targetId in $idsToDelete) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Id -eq $targetId) {
            $shp.Delete()
            break
        }
    }
}
